$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: use the same column names as the Jira export (add spaces) ---
$ws.Range("A1").Value = "Sprint Name"
$ws.Range("B1").Value = "Start Date"
$ws.Range("C1").Value = "End Date"
$ws.Range("D1").Value = "Capacity Forecast"
$ws.Range("E1").Value = "Effort Forecast"
$ws.Range("F1").Value = "Capacity Done"
$ws.Range("G1").Value = "Effort Done"

# --- Remove the helper "Forecast"/"Done" formula columns (I:L) ---
$ws.Range("I1:L9").Clear()

# --- Data updates ---
# Row 5 (Sprint 4): CapacityDone / EffortDone
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 2

# Row 6 (Sprint 5): End Date, CapacityForecast, EffortForecast
$ws.Range("C6").Value = 41763
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 5

# Row 7 (Sprint 6): Start Date, End Date, CapacityForecast
$ws.Range("B7").Value = 41764
$ws.Range("C7").Value = 41770
$ws.Range("D7").Value = 12

# Row 8 (Sprint 7): Start Date, End Date, CapacityForecast
$ws.Range("B8").Value = 41771
$ws.Range("C8").Value = 41777
$ws.Range("D8").Value = 12

# Row 9 (Sprint 8): Start Date, End Date, CapacityForecast
$ws.Range("B9").Value = 41778
$ws.Range("C9").Value = 41784
$ws.Range("D9").Value = 12

# --- View: selection moved away from the removed helper-formula area ---
$ws.Range("H1:L6").Select()
